# 10Th - MB for single stock and added new group
#
# The weekly MarketBeat rank sheet rolls forward: two brand-new date columns
# (Jun_26, Jun_26, Jun_27 -> really 3 new columns B:D) are inserted in front
# of the existing date columns, which shift right, and two new brokerage
# rows (Benchmark, Evercore ISI) are appended at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new columns before column B, shifting the old B:E -> E:H ---
$ws.Range("B1:D27").Insert(-4161)   # xlShiftToRight

# Keep the column widths consistent (~8 chars) across the newly inserted
# and the pre-existing date columns, same as the rest of the table.
$ws.Columns("C:H").ColumnWidth = 7.14

# --- New header values for the freshly inserted date columns ---
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# --- Fill the new columns for every existing brokerage row with "UN" ---
$ws.Range("B2:D27").Value = "UN"

# --- Append two new brokerage rows at the bottom ---
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28:D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29:D29").Value = "UN"
